$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7, column B: reuse the shaded/centered label style already used by B5 & B6
$ws.Range("B5").Copy()
$ws.Range("B7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B7").Value = "drim"

# C7: currency amount with bold-ish (10pt) font, horizontally centered
$ws.Range("C7").Value = 9.48
$ws.Range("C7").NumberFormat = "#,##0.00\ ""€"";[Red]\-#,##0.00\ ""€"""
$ws.Range("C7").HorizontalAlignment = -4108
$ws.Range("C7").VerticalAlignment = -4107
$ws.Range("C7").Font.Size = 10

# D7: currency amount, default font, horizontally centered
$ws.Range("D7").Value = 9.48
$ws.Range("D7").NumberFormat = "#,##0.00\ ""€"";[Red]\-#,##0.00\ ""€"""
$ws.Range("D7").HorizontalAlignment = -4108
$ws.Range("D7").VerticalAlignment = -4107

# Update selected cell as last user action
$ws.Range("K15").Select() | Out-Null
